# Update the "Cereza" (cherry) price-report worksheet:
#  - rows 58-60 become new "Early Burlat" entries
#  - rows 61-70 are overwritten with the data that used to sit 3 rows above
#    them (58-67), i.e. the block of rows is effectively shifted down by
#    three positions to make room for the new entries
#  - three brand-new rows (72-74) are appended at the bottom, reusing what
#    used to be in rows 69-71
#  - the very last row (74) keeps the original row-71 values untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that stay constant for every one of these "Cereza" rows.
$colA = 5
$colB = "Macroferia Regional de Talca"
$colC = "Maule"
$colE = 7
$colF = "Fruta"
$colG = 100103
$colH = "Frutos de hueso (carozo)"
$colI = 100103001
$colJ = "Cereza"
$colR = "Provincia de Curicó"

# Target state for every row from 58 to 74 (inclusive), keyed by row number.
$rows = @{
    58 = @{ D = 44511; K = "Early Burlat"; L = "Primera"; M = 50;  N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 5 kilos";            S = 3000; T = 5 }
    59 = @{ D = 44511; K = "Early Burlat"; L = "Primera"; M = 100; N = 20000; O = 20000; P = 20000; Q = "`$/bandeja 7 kilos";            S = 2857; T = 7 }
    60 = @{ D = 44511; K = "Early Burlat"; L = "Primera"; M = 20;  N = 3000;  O = 3000;  P = 3000;  Q = "`$/kilo (en caja de 15 kilos)"; S = 3000; T = 1 }
    61 = @{ D = 44211; K = "Bing";         L = "Primera"; M = 120; N = 5000;  O = 5000;  P = 5000;  Q = "`$/bandeja 10 kilos";           S = 500;  T = 10 }
    62 = @{ D = 44211; K = "Sweet Heart";  L = "Primera"; M = 140; N = 6000;  O = 6000;  P = 6000;  Q = "`$/bandeja 10 kilos";           S = 600;  T = 10 }
    63 = @{ D = 44186; K = "Bing";         L = "Primera"; M = 200; N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 10 kilos";              S = 800;  T = 10 }
    64 = @{ D = 44186; K = "Lapins";       L = "Primera"; M = 200; N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 10 kilos";              S = 800;  T = 10 }
    65 = @{ D = 44168; K = "Lapins";       L = "Primera"; M = 350; N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 10 kilos";           S = 1000; T = 10 }
    66 = @{ D = 44168; K = "Rainier";      L = "Primera"; M = 150; N = 13000; O = 13000; P = 13000; Q = "`$/caja 10 kilos";              S = 1300; T = 10 }
    67 = @{ D = 44168; K = "Santina";      L = "Primera"; M = 160; N = 10000; O = 10000; P = 10000; Q = "`$/caja 10 kilos";              S = 1000; T = 10 }
    68 = @{ D = 44175; K = "Rainier";      L = "Primera"; M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 10 kilos";              S = 1200; T = 10 }
    69 = @{ D = 44175; K = "Santina";      L = "Primera"; M = 100; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 10 kilos";              S = 900;  T = 10 }
    70 = @{ D = 44175; K = "Santina";      L = "Segunda"; M = 60;  N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 10 kilos";              S = 800;  T = 10 }
    71 = @{ D = 44181; K = "Lapins";       L = "Primera"; M = 100; N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 10 kilos";              S = 800;  T = 10 }
    72 = @{ D = 44181; K = "Santina";      L = "Primera"; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 10 kilos";              S = 900;  T = 10 }
    73 = @{ D = 44217; K = "Bing";         L = "Primera"; M = 150; N = 6000;  O = 7000;  P = 6667;  Q = "`$/bandeja 10 kilos";           S = 667;  T = 10 }
    74 = @{ D = 44179; K = "Lapins";       L = "Primera"; M = 450; N = 9000;  O = 10000; P = 9444;  Q = "`$/bandeja 10 kilos";           S = 944;  T = 10 }
}

foreach ($r in 58..74) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC

    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ

    $ws.Cells.Item($r, 11).Value = $vals.K
    $ws.Cells.Item($r, 12).Value = $vals.L
    $ws.Cells.Item($r, 13).Value = $vals.M
    $ws.Cells.Item($r, 14).Value = $vals.N
    $ws.Cells.Item($r, 15).Value = $vals.O
    $ws.Cells.Item($r, 16).Value = $vals.P
    $ws.Cells.Item($r, 17).Value = $vals.Q

    $ws.Cells.Item($r, 18).Value = $colR

    $ws.Cells.Item($r, 19).Value = $vals.S
    $ws.Cells.Item($r, 20).Value = $vals.T
}
